$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 3 wrapped up - fill in the column D results that just came in
$ws.Range("D3").Value = 3
$ws.Range("D5").Value = 1
$ws.Range("D7").Value = "1.  Government manage all opportunities"

# The new, longer "Features Accepted" entry makes row 7 wrap onto as many
# lines as row 6 (same text length/style), so it grows to match that height.
$ws.Rows(7).RowHeight = 76.5

# Column D widened slightly to better fit the newly-entered data.
$ws.Columns(4).ColumnWidth = 13.451822916666666

# Leave the selection where data entry ended up: D6:D7, active cell D6.
$excel.Goto($ws.Range("D6:D7"))
